$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "Índice"
$ws.Range("B1").Value = "Distancia"
$ws.Range("C1").Value = "max"
$ws.Range("D1").Value = "min"
$ws.Range("E1").Value = "Tempo"

# Data rows
$data = @(
    @(0, 2031, 2031, 2031, 0.01132214069366455),
    @(1, 2296, 2296, 2296, 0.01333446502685547),
    @(2, 1822, 1822, 1822, 0.01333240667978922),
    @(3, 2829, 2829, 2829, 0.01301171779632568),
    @(4, 2187, 2187, 2187, 0.01325689951578776),
    @(5, 3125, 3125, 3125, 0.01344212690989176),
    @(6, 2628, 2628, 2628, 0.01284503936767578),
    @(7, 2734, 2734, 2734, 0.01339208285013835),
    @(8, 2886, 2886, 2886, 0.01351958910624186),
    @(9, 2917, 2917, 2917, 0.01355818112691243)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $rowData = $data[$i]
    $ws.Cells.Item($row, 1).Value = $rowData[0]
    $ws.Cells.Item($row, 2).Value = $rowData[1]
    $ws.Cells.Item($row, 3).Value = $rowData[2]
    $ws.Cells.Item($row, 4).Value = $rowData[3]
    $ws.Cells.Item($row, 5).Value = $rowData[4]
}
